# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt - Perejil"
# right after the current row 139 (i.e. as the new row 140), pushing the
# existing rows 140:173 down to 141:174.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 140:173 down to 141:174, leaving a blank row 140 behind.
$ws.Rows(140).Insert()

# Populate the new row 140 with this week's observation.
$ws.Cells.Item(140, 1).Value = 4
$ws.Cells.Item(140, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(140, 3).Value = "Los Lagos"
$ws.Cells.Item(140, 4).Value = 44511
$ws.Cells.Item(140, 5).Value = 10
$ws.Cells.Item(140, 6).Value = 100112044
$ws.Cells.Item(140, 7).Value = "Perejil"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 60
$ws.Cells.Item(140, 11).Value = 5000
$ws.Cells.Item(140, 12).Value = 5000
$ws.Cells.Item(140, 13).Value = 5000
$ws.Cells.Item(140, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(140, 15).Value = "Región Metropolitana"
$ws.Cells.Item(140, 16).Value = 1667
$ws.Cells.Item(140, 17).Value = 3
$ws.Cells.Item(140, 18).Value = "Hortaliza"
